# Generate Report for handoff
# A new source file (02dafc58-253b-4d5e-9f9e-e7a9112e8baf.md) shows up in the
# localization report with a "Handoff transform failed" status, and the
# already-tracked file got a new guid (60923cb0-ac7c-4e4a-9902-e7ebc6f1600a.md)
# together with a refreshed handoff package / timestamps.

$wb = $excel.ActiveWorkbook

# Cornflower blue (FF6495ED) expressed as the BGR integer Excel's Font.Color expects.
$hyperlinkColor = 15570276

function Style-AsText {
    param($ws, $addr)
    # Plain body text cell (style 0 in the original workbook) - nothing to do,
    # General number format / default font is already what we want.
}

function Style-AsHyperlink {
    param($ws, $addr)
    $r = $ws.Range($addr)
    $r.Font.Underline = 2
    $r.Font.Color = $hyperlinkColor
    $r.Font.Name = "Calibri"
}

function Style-AsDate {
    param($ws, $addr)
    $ws.Range($addr).NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

function Add-Link {
    param($ws, $addr, $url, $display)
    $ws.Hyperlinks.Add($ws.Range($addr), $url, "", "", $display) | Out-Null
    Style-AsHyperlink $ws $addr
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Hyperlinks.Delete()

$ws1.Range("A2").Value = "60923cb0-ac7c-4e4a-9902-e7ebc6f1600a.md"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

$ws1.Range("A3").Value = "02dafc58-253b-4d5e-9f9e-e7a9112e8baf.md"
$ws1.Range("B3").Value = "Handoff transform failed"
$ws1.Range("C3").Value = "Handoff transform failed"

$ws1.Range("A4").Value = ".localization-config"
$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

Add-Link $ws1 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/f33a6c89301881f2e6d70ed4d8b7ba55ca2bda2b/e2e/60923cb0-ac7c-4e4a-9902-e7ebc6f1600a.md" "60923cb0-ac7c-4e4a-9902-e7ebc6f1600a.md"
Add-Link $ws1 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/f33a6c89301881f2e6d70ed4d8b7ba55ca2bda2b/e2e/02dafc58-253b-4d5e-9f9e-e7a9112e8baf.md" "02dafc58-253b-4d5e-9f9e-e7a9112e8baf.md"
Add-Link $ws1 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/f33a6c89301881f2e6d70ed4d8b7ba55ca2bda2b/.localization-config" ".localization-config"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()

$ws2.Range("A2").Value = "60923cb0-ac7c-4e4a-9902-e7ebc6f1600a.md"
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("C2").Value = "60923cb0-ac7c-4e4a-9902-e7ebc6f1600a.828a692fff8a8a081b8584214f38fc9f2626e92.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-01-14 05:31:36"
$ws2.Range("G2").Value = "0001-01-01 00:00:00"
$ws2.Range("H2").Value = "Include"
Style-AsDate $ws2 "D2"

$ws2.Range("A3").Value = "02dafc58-253b-4d5e-9f9e-e7a9112e8baf.md"
$ws2.Range("B3").Value = "Handoff transform failed"
$ws2.Range("D3").Value = "0001-01-01 00:00:00"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Ignored"
Style-AsDate $ws2 "D3"

$ws2.Range("A4").Value = ".localization-config"
$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("D4").Value = "0001-01-01 00:00:00"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Ignored"
Style-AsDate $ws2 "D4"

Add-Link $ws2 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/f33a6c89301881f2e6d70ed4d8b7ba55ca2bda2b/e2e/60923cb0-ac7c-4e4a-9902-e7ebc6f1600a.md" "60923cb0-ac7c-4e4a-9902-e7ebc6f1600a.md"
Add-Link $ws2 "C2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d7693261d974fb786be6d47b219cd97c4352b6a9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/60923cb0-ac7c-4e4a-9902-e7ebc6f1600a.828a692fff8a8a081b8584214f38fc9f2626e92.zh-cn.xlf" "60923cb0-ac7c-4e4a-9902-e7ebc6f1600a.828a692fff8a8a081b8584214f38fc9f2626e92.zh-cn.xlf"
Add-Link $ws2 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/f33a6c89301881f2e6d70ed4d8b7ba55ca2bda2b/e2e/02dafc58-253b-4d5e-9f9e-e7a9112e8baf.md" "02dafc58-253b-4d5e-9f9e-e7a9112e8baf.md"
Add-Link $ws2 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/f33a6c89301881f2e6d70ed4d8b7ba55ca2bda2b/.localization-config" ".localization-config"

# re-apply the date number format on the linked cell; Hyperlinks.Add can reset direct formatting
Style-AsDate $ws2 "D2"
Style-AsDate $ws2 "D3"
Style-AsDate $ws2 "D4"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()

$ws3.Range("A2").Value = "60923cb0-ac7c-4e4a-9902-e7ebc6f1600a.md"
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("C2").Value = "60923cb0-ac7c-4e4a-9902-e7ebc6f1600a.828a692fff8a8a081b8584214f38fc9f2626e92.de-de.xlf"
$ws3.Range("D2").Value = "2016-01-14 05:31:59"
$ws3.Range("G2").Value = "0001-01-01 00:00:00"
$ws3.Range("H2").Value = "Include"
Style-AsDate $ws3 "D2"

$ws3.Range("A3").Value = "02dafc58-253b-4d5e-9f9e-e7a9112e8baf.md"
$ws3.Range("B3").Value = "Handoff transform failed"
$ws3.Range("D3").Value = "0001-01-01 00:00:00"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Ignored"
Style-AsDate $ws3 "D3"

$ws3.Range("A4").Value = ".localization-config"
$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("D4").Value = "0001-01-01 00:00:00"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Ignored"
Style-AsDate $ws3 "D4"

Add-Link $ws3 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/f33a6c89301881f2e6d70ed4d8b7ba55ca2bda2b/e2e/60923cb0-ac7c-4e4a-9902-e7ebc6f1600a.md" "60923cb0-ac7c-4e4a-9902-e7ebc6f1600a.md"
Add-Link $ws3 "C2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8acb6804190dae2fd00a4e0ccc0be2b96b7a834b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/60923cb0-ac7c-4e4a-9902-e7ebc6f1600a.828a692fff8a8a081b8584214f38fc9f2626e92.de-de.xlf" "60923cb0-ac7c-4e4a-9902-e7ebc6f1600a.828a692fff8a8a081b8584214f38fc9f2626e92.de-de.xlf"
Add-Link $ws3 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/f33a6c89301881f2e6d70ed4d8b7ba55ca2bda2b/e2e/02dafc58-253b-4d5e-9f9e-e7a9112e8baf.md" "02dafc58-253b-4d5e-9f9e-e7a9112e8baf.md"
Add-Link $ws3 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/f33a6c89301881f2e6d70ed4d8b7ba55ca2bda2b/.localization-config" ".localization-config"

Style-AsDate $ws3 "D2"
Style-AsDate $ws3 "D3"
Style-AsDate $ws3 "D4"

Write-Output "Applied handoff report update"
